$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.810.34"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.906.07"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -4.19%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.80"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.22"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -5.99%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.502"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.908.07"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.445"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -4.01%  "
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.38"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -6.60%  "
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.393.43"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.772.94"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.70"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -5.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.911.41"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -4.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "427.76"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -5.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.52"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("E23").Value = "  -5.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.49"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.80"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.89"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -4.32%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.16"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -4.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.44"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.105"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0873"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("E36").Value = "  -3.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.00"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -5.85%  "
$ws.Range("E39").Value = "  -3.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.42"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  -5.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.58"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -6.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.295"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.37"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -5.79%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0350"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "377.28"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.687.18"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.16"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.34"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  -2.80%  "
